$wb = $excel.ActiveWorkbook

# New address text for 李宁运动中心 location (replaces old 机场路1399号... address)
$newAddress = "广龙路中油BP(白云万顺达南加油站)北侧约260米 李宁运动中心"

# ---- Sheet: 展览 ----
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 915
$ws.Range("F4").Value = 807
$ws.Range("D5").Value = $newAddress
$ws.Range("F6").Value = 463
$ws.Range("F7").Value = 712
$ws.Range("F8").Value = 166
$ws.Range("F9").Value = 1314
$ws.Range("F10").Value = 732
$ws.Range("F11").Value = 423
$ws.Range("D12").Value = $newAddress
$ws.Range("F13").Value = 189
$ws.Range("F14").Value = 55
$ws.Range("F15").Value = 1164
$ws.Range("F21").Value = 601
$ws.Range("F23").Value = 665
$ws.Range("F25").Value = 1079

# ---- Sheet: 演出 ----
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F2").Value = 348
$ws.Range("F7").Value = 256

# ---- Sheet: 本地生活 ----
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 381

# ---- Sheet: 全部类型 ----
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 381
$ws.Range("F3").Value = 348
$ws.Range("F4").Value = 915
$ws.Range("F6").Value = 807
$ws.Range("D7").Value = $newAddress
$ws.Range("F8").Value = 463
$ws.Range("F9").Value = 463
$ws.Range("F10").Value = 712
$ws.Range("F11").Value = 166
$ws.Range("F12").Value = 1314
$ws.Range("F13").Value = 732
$ws.Range("F16").Value = 423
$ws.Range("D17").Value = $newAddress
$ws.Range("F19").Value = 189
$ws.Range("F20").Value = 55
$ws.Range("F21").Value = 1164
$ws.Range("F28").Value = 256
$ws.Range("F30").Value = 601
$ws.Range("F36").Value = 665
$ws.Range("F38").Value = 1079
